# Auto-generated edit script applying numeric cell updates per diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15 (ALC) - hunk 0
$ws.Cells.Item(15, 8).Value = 572.3111
$ws.Cells.Item(15, 9).Value = 572.3111
$ws.Cells.Item(15, 11).Value = 1716.9333
$ws.Cells.Item(15, 13).Value = -1547.9333

# Row 32 (ALC) - hunk 1
$ws.Cells.Item(32, 8).Value = 1700.25
$ws.Cells.Item(32, 10).Value = 0
$ws.Cells.Item(32, 12).Value = 0
$ws.Cells.Item(32, 14).ClearContents()

# Row 100 (ALC) - hunk 2
$ws.Cells.Item(100, 8).Value = 2946.5334
$ws.Cells.Item(100, 9).Value = 2854
$ws.Cells.Item(100, 10).Value = 3316.6667
$ws.Cells.Item(100, 11).Value = 2854
$ws.Cells.Item(100, 12).Value = 3316.6667
$ws.Cells.Item(100, 13).Value = -2313
$ws.Cells.Item(100, 14).Value = -4398.6667

# Row 103 (ALC) - hunk 3
$ws.Cells.Item(103, 8).Value = 707.5
$ws.Cells.Item(103, 10).Value = 707.5
$ws.Cells.Item(103, 12).Value = 2122.5
$ws.Cells.Item(103, 14).Value = -3294.5

# Row 116 (ALC) - hunk 4
$ws.Cells.Item(116, 8).Value = 4770.8335
$ws.Cells.Item(116, 9).Value = 2725
$ws.Cells.Item(116, 10).Value = 15000
$ws.Cells.Item(116, 11).Value = 2725
$ws.Cells.Item(116, 12).Value = 15000
$ws.Cells.Item(116, 13).Value = 717
$ws.Cells.Item(116, 14).Value = -21884

# Row 138 (ALC) - hunk 5
$ws.Cells.Item(138, 8).Value = 5320.8604
$ws.Cells.Item(138, 9).Value = 5031.0713
$ws.Cells.Item(138, 10).Value = 5460.759
$ws.Cells.Item(138, 11).Value = 15093.2139
$ws.Cells.Item(138, 12).Value = 16382.277
$ws.Cells.Item(138, 13).Value = -9953.213899999999
$ws.Cells.Item(138, 14).Value = -26662.277

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (ARM) - hunk 6
$ws.Cells.Item(32, 8).Value = 18063.072
$ws.Cells.Item(32, 10).Value = 43299.332
$ws.Cells.Item(32, 12).Value = 43299.332
$ws.Cells.Item(32, 14).Value = -43873.332

# Row 35 (ARM) - hunk 7
$ws.Cells.Item(35, 8).Value = 0
$ws.Cells.Item(35, 9).Value = 0
$ws.Cells.Item(35, 11).Value = 0
$ws.Cells.Item(35, 13).ClearContents()

# Row 61 (ARM) - hunk 8
$ws.Cells.Item(61, 8).Value = 6270.8423
$ws.Cells.Item(61, 9).Value = 6527.9375
$ws.Cells.Item(61, 11).Value = 6527.9375
$ws.Cells.Item(61, 13).Value = -6315.9375

# Row 63 (ARM) - hunk 9
$ws.Cells.Item(63, 8).Value = 1066.6666
$ws.Cells.Item(63, 9).Value = 1066.6666
$ws.Cells.Item(63, 10).Value = 0
$ws.Cells.Item(63, 11).Value = 1066.6666
$ws.Cells.Item(63, 12).Value = 0
$ws.Cells.Item(63, 13).Value = -380.6666
$ws.Cells.Item(63, 14).ClearContents()

# Row 66 (ARM) - hunk 10
$ws.Cells.Item(66, 8).Value = 1066.6666
$ws.Cells.Item(66, 9).Value = 1066.6666
$ws.Cells.Item(66, 10).Value = 0
$ws.Cells.Item(66, 11).Value = 5333.333000000001
$ws.Cells.Item(66, 12).Value = 0
$ws.Cells.Item(66, 13).Value = -1901.333000000001
$ws.Cells.Item(66, 14).ClearContents()

# Row 122 (ARM) - hunk 11
$ws.Cells.Item(122, 8).Value = 4674.1577
$ws.Cells.Item(122, 9).Value = 4433.1763
$ws.Cells.Item(122, 11).Value = 13299.5289
$ws.Cells.Item(122, 13).Value = -10849.5289

# Row 132 (ARM) - hunk 12
$ws.Cells.Item(132, 8).Value = 5699.4546
$ws.Cells.Item(132, 9).Value = 5699.4546
$ws.Cells.Item(132, 11).Value = 17098.3638
$ws.Cells.Item(132, 13).Value = -14568.3638

# Row 136 (ARM) - hunk 13
$ws.Cells.Item(136, 8).Value = 6270.8423
$ws.Cells.Item(136, 9).Value = 6527.9375
$ws.Cells.Item(136, 11).Value = 19583.8125
$ws.Cells.Item(136, 13).Value = -17033.8125

$ws = $wb.Worksheets.Item("BSM")
# Row 94 (BSM) - hunk 14
$ws.Cells.Item(94, 8).Value = 686.4375
$ws.Cells.Item(94, 9).Value = 652.2
$ws.Cells.Item(94, 11).Value = 652.2
$ws.Cells.Item(94, 13).Value = -201.2

# Row 105 (BSM) - hunk 15
$ws.Cells.Item(105, 8).Value = 1378.55
$ws.Cells.Item(105, 9).Value = 1223.5
$ws.Cells.Item(105, 11).Value = 1223.5
$ws.Cells.Item(105, 13).Value = 523.5

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (CRP) - hunk 16
$ws.Cells.Item(31, 8).Value = 2888.6667
$ws.Cells.Item(31, 10).Value = 3166.6667
$ws.Cells.Item(31, 12).Value = 3166.6667
$ws.Cells.Item(31, 14).Value = -3756.6667

# Row 34 (CRP) - hunk 17
$ws.Cells.Item(34, 8).Value = 2888.6667
$ws.Cells.Item(34, 10).Value = 3166.6667
$ws.Cells.Item(34, 12).Value = 3166.6667
$ws.Cells.Item(34, 14).Value = -3570.6667

# Row 58 (CRP) - hunk 18
$ws.Cells.Item(58, 8).Value = 2258.6191
$ws.Cells.Item(58, 9).Value = 2002.1177
$ws.Cells.Item(58, 10).Value = 3348.75
$ws.Cells.Item(58, 11).Value = 2002.1177
$ws.Cells.Item(58, 12).Value = 3348.75
$ws.Cells.Item(58, 13).Value = -1799.1177
$ws.Cells.Item(58, 14).Value = -3754.75

# Row 62 (CRP) - hunk 19
$ws.Cells.Item(62, 8).Value = 3546.25
$ws.Cells.Item(62, 9).Value = 2866.6667
$ws.Cells.Item(62, 10).Value = 5585
$ws.Cells.Item(62, 11).Value = 2866.6667
$ws.Cells.Item(62, 12).Value = 5585
$ws.Cells.Item(62, 13).Value = -2242.6667
$ws.Cells.Item(62, 14).Value = -6833

# Row 65 (CRP) - hunk 20
$ws.Cells.Item(65, 8).Value = 3546.25
$ws.Cells.Item(65, 9).Value = 2866.6667
$ws.Cells.Item(65, 10).Value = 5585
$ws.Cells.Item(65, 11).Value = 14333.3335
$ws.Cells.Item(65, 12).Value = 27925
$ws.Cells.Item(65, 13).Value = -11213.3335
$ws.Cells.Item(65, 14).Value = -34165

# Row 94 (CRP) - hunk 21
$ws.Cells.Item(94, 8).Value = 0
$ws.Cells.Item(94, 9).Value = 0
$ws.Cells.Item(94, 11).Value = 0
$ws.Cells.Item(94, 13).ClearContents()

# Row 132 (CRP) - hunk 22
$ws.Cells.Item(132, 8).Value = 14295886
$ws.Cells.Item(132, 10).Value = 2999
$ws.Cells.Item(132, 12).Value = 8997
$ws.Cells.Item(132, 14).Value = -14057

# Row 136 (CRP) - hunk 23
$ws.Cells.Item(136, 8).Value = 2258.6191
$ws.Cells.Item(136, 9).Value = 2002.1177
$ws.Cells.Item(136, 10).Value = 3348.75
$ws.Cells.Item(136, 11).Value = 6006.3531
$ws.Cells.Item(136, 12).Value = 10046.25
$ws.Cells.Item(136, 13).Value = -3456.3531
$ws.Cells.Item(136, 14).Value = -15146.25

$ws = $wb.Worksheets.Item("CUL")
# Row 5 (CUL) - hunk 24
$ws.Cells.Item(5, 8).Value = 69505.8
$ws.Cells.Item(5, 9).Value = 3160.1538
$ws.Cells.Item(5, 10).Value = 500752.5
$ws.Cells.Item(5, 11).Value = 9480.4614
$ws.Cells.Item(5, 12).Value = 1502257.5
$ws.Cells.Item(5, 13).Value = -9368.4614
$ws.Cells.Item(5, 14).Value = -1502481.5

# Row 44 (CUL) - hunk 25
$ws.Cells.Item(44, 8).Value = 111331.22
$ws.Cells.Item(44, 9).Value = 166706.83
$ws.Cells.Item(44, 11).Value = 500120.49
$ws.Cells.Item(44, 13).Value = -499722.49

# Row 46 (CUL) - hunk 26
$ws.Cells.Item(46, 8).Value = 0
$ws.Cells.Item(46, 9).Value = 0
$ws.Cells.Item(46, 10).Value = 0
$ws.Cells.Item(46, 11).Value = 0
$ws.Cells.Item(46, 12).Value = 0
$ws.Cells.Item(46, 13).ClearContents()
$ws.Cells.Item(46, 14).ClearContents()

# Row 55 (CUL) - hunk 27
$ws.Cells.Item(55, 8).Value = 3102.6
$ws.Cells.Item(55, 9).Value = 401.5
$ws.Cells.Item(55, 10).Value = 4903.3335
$ws.Cells.Item(55, 11).Value = 1204.5
$ws.Cells.Item(55, 12).Value = 14710.0005
$ws.Cells.Item(55, 13).Value = -1027.5
$ws.Cells.Item(55, 14).Value = -15064.0005

# Row 69 (CUL) - hunk 28
$ws.Cells.Item(69, 8).Value = 2662.093
$ws.Cells.Item(69, 9).Value = 3096.6
$ws.Cells.Item(69, 11).Value = 9289.799999999999
$ws.Cells.Item(69, 13).Value = -8478.799999999999

# Row 72 (CUL) - hunk 29
$ws.Cells.Item(72, 8).Value = 2662.093
$ws.Cells.Item(72, 9).Value = 3096.6
$ws.Cells.Item(72, 11).Value = 27869.4
$ws.Cells.Item(72, 13).Value = -23813.4

# Row 117 (CUL) - hunk 30
$ws.Cells.Item(117, 8).Value = 995.4
$ws.Cells.Item(117, 10).Value = 1177
$ws.Cells.Item(117, 12).Value = 3531
$ws.Cells.Item(117, 14).Value = -10415

# Row 135 (CUL) - hunk 31
$ws.Cells.Item(135, 8).Value = 69505.8
$ws.Cells.Item(135, 9).Value = 3160.1538
$ws.Cells.Item(135, 10).Value = 500752.5
$ws.Cells.Item(135, 11).Value = 28441.3842
$ws.Cells.Item(135, 12).Value = 4506772.5
$ws.Cells.Item(135, 13).Value = -25906.3842
$ws.Cells.Item(135, 14).Value = -4511842.5

$ws = $wb.Worksheets.Item("GSM")
# Row 2 (GSM) - hunk 32
$ws.Cells.Item(2, 8).Value = 766.4211
$ws.Cells.Item(2, 9).Value = 984.75
$ws.Cells.Item(2, 11).Value = 984.75
$ws.Cells.Item(2, 13).Value = -871.75

# Row 125 (GSM) - hunk 33
$ws.Cells.Item(125, 8).Value = 0
$ws.Cells.Item(125, 10).Value = 0
$ws.Cells.Item(125, 12).Value = 0
$ws.Cells.Item(125, 14).ClearContents()

# Row 132 (GSM) - hunk 34
$ws.Cells.Item(132, 8).Value = 27780202
$ws.Cells.Item(132, 9).Value = 2658.8
$ws.Cells.Item(132, 10).Value = 166667920
$ws.Cells.Item(132, 11).Value = 7976.400000000001
$ws.Cells.Item(132, 12).Value = 500003760
$ws.Cells.Item(132, 13).Value = -5446.400000000001
$ws.Cells.Item(132, 14).Value = -500008820

$ws = $wb.Worksheets.Item("LTW")
# Row 46 (LTW) - hunk 35
$ws.Cells.Item(46, 8).Value = 1244.5
$ws.Cells.Item(46, 9).Value = 1190
$ws.Cells.Item(46, 11).Value = 1190
$ws.Cells.Item(46, 13).Value = -1002

# Row 132 (LTW) - hunk 36
$ws.Cells.Item(132, 8).Value = 7875.8
$ws.Cells.Item(132, 9).Value = 5154.6
$ws.Cells.Item(132, 10).Value = 13318.2
$ws.Cells.Item(132, 11).Value = 15463.8
$ws.Cells.Item(132, 12).Value = 39954.60000000001
$ws.Cells.Item(132, 13).Value = -12933.8
$ws.Cells.Item(132, 14).Value = -45014.60000000001

# Row 136 (LTW) - hunk 37
$ws.Cells.Item(136, 8).Value = 2188.842
$ws.Cells.Item(136, 9).Value = 2158.1177
$ws.Cells.Item(136, 11).Value = 6474.353099999999
$ws.Cells.Item(136, 13).Value = -3924.353099999999

$ws = $wb.Worksheets.Item("WVR")
# Row 132 (WVR) - hunk 38
$ws.Cells.Item(132, 8).Value = 1000000000
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 11).Value = 0
$ws.Cells.Item(132, 13).ClearContents()
